$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16504778589745333"
$ws1.Range("B2").Value = "go_stims-16504778589325314.csv"
$ws1.Range("B3").Value = "GNG_stims-1650477858958532.csv"
$ws1.Range("B4").Value = "go_stims-16504778589595306.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778589735322.csv"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16504778606946154"
$ws2.Range("B2").Value = "ZB-match_2-16504778593735297.csv"
$ws2.Range("B3").Value = "ZB-match_1-165047785926653.csv"
$ws2.Range("B4").Value = "TB-16504778605355613.csv"
$ws2.Range("B5").Value = "OB-16504778600185273.csv"
$ws2.Range("B6").Value = "TB-16504778602245622.csv"
$ws2.Range("B7").Value = "OB-16504778597655625.csv"
$ws2.Range("B8").Value = "ZB-match_1-16504778589925287.csv"
$ws2.Range("B9").Value = "OB-16504778599055638.csv"
$ws2.Range("B10").Value = "TB-16504778606805594.csv"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16504778606955323"
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16504778607435274"
$ws4.Range("B2").Value = "MM_stims-16504778607105646.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778606975267.csv"
$ws4.Range("B4").Value = "MM_stims-16504778607265291.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778607105646.csv"
$ws4.Range("B6").Value = "MM_stims-16504778607425652.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778607265291.csv"

$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16504778608065317"
$ws5.Range("B2").Value = "vSAT_stims-16504778607745283.csv"
$ws5.Range("B3").Value = "SAT_stims-16504778607465298.csv"
$ws5.Range("B4").Value = "SAT_stims-16504778607585652.csv"
$ws5.Range("B5").Value = "vSAT_stims-165047786079053.csv"
